$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10151.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 10151.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 10151.5
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -11119.5

$ws.Range("H64").Value = 2887.3333
$ws.Range("I64").Value = 2800
$ws.Range("J64").Value = 2898.25
$ws.Range("K64").Value = 2800
$ws.Range("L64").Value = 2898.25
$ws.Range("M64").Value = -2552
$ws.Range("N64").Value = -3394.25

$ws.Range("H67").Value = 2887.3333
$ws.Range("I67").Value = 2800
$ws.Range("J67").Value = 2898.25
$ws.Range("K67").Value = 2800
$ws.Range("L67").Value = 2898.25
$ws.Range("M67").Value = -1942
$ws.Range("N67").Value = -4614.25

$ws.Range("H74").Value = 4078.2666
$ws.Range("I74").Value = 4097.222
$ws.Range("J74").Value = 4049.8333
$ws.Range("K74").Value = 4097.222
$ws.Range("L74").Value = 4049.8333
$ws.Range("M74").Value = -3161.222
$ws.Range("N74").Value = -5921.8333

$ws.Range("H77").Value = 4078.2666
$ws.Range("I77").Value = 4097.222
$ws.Range("J77").Value = 4049.8333
$ws.Range("K77").Value = 20486.11
$ws.Range("L77").Value = 20249.1665
$ws.Range("M77").Value = -15806.11
$ws.Range("N77").Value = -29609.1665

$ws.Range("H112").Value = 3847.7036
$ws.Range("J112").Value = 3838
$ws.Range("L112").Value = 11514
$ws.Range("N112").Value = -13730

$ws.Range("H127").Value = 2810.9736
$ws.Range("I127").Value = 1016.25
$ws.Range("J127").Value = 3289.5667
$ws.Range("K127").Value = 3048.75
$ws.Range("L127").Value = 9868.7001
$ws.Range("M127").Value = 1911.25
$ws.Range("N127").Value = -19788.7001

$ws.Range("H131").Value = 2525.6667
$ws.Range("I131").Value = 937.0909
$ws.Range("K131").Value = 2811.2727
$ws.Range("M131").Value = 2228.7273

$ws.Range("H132").Value = 4660.6665
$ws.Range("I132").Value = 1458.5172
$ws.Range("J132").Value = 17926.715
$ws.Range("K132").Value = 4375.5516
$ws.Range("L132").Value = 53780.145
$ws.Range("M132").Value = -1845.5516
$ws.Range("N132").Value = -58840.145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15149.217
$ws.Range("I32").Value = 14108.272
$ws.Range("K32").Value = 14108.272
$ws.Range("M32").Value = -13821.272

$ws.Range("H45").Value = 2765.9092
$ws.Range("I45").Value = 1604.5
$ws.Range("K45").Value = 1604.5
$ws.Range("M45").Value = -1227.5

$ws.Range("H122").Value = 5988.4
$ws.Range("I122").Value = 5431.5557
$ws.Range("K122").Value = 16294.6671
$ws.Range("M122").Value = -13844.6671

$ws.Range("H132").Value = 2251.457
$ws.Range("I132").Value = 2251.457
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6754.370999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4224.370999999999
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2260.45
$ws.Range("I107").Value = 2174
$ws.Range("J107").Value = 2440
$ws.Range("K107").Value = 2174
$ws.Range("L107").Value = 2440
$ws.Range("M107").Value = -254
$ws.Range("N107").Value = -6280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14708678
$ws.Range("I31").Value = 15875424
$ws.Range("J31").Value = 7675
$ws.Range("K31").Value = 15875424
$ws.Range("L31").Value = 7675
$ws.Range("M31").Value = -15875129
$ws.Range("N31").Value = -8265

$ws.Range("H34").Value = 14708678
$ws.Range("I34").Value = 15875424
$ws.Range("J34").Value = 7675
$ws.Range("K34").Value = 15875424
$ws.Range("L34").Value = 7675
$ws.Range("M34").Value = -15875222
$ws.Range("N34").Value = -8079

$ws.Range("H99").Value = 20001.2
$ws.Range("I99").Value = 110012
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 110012
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -108514
$ws.Range("N99").Value = -12996

$ws.Range("H102").Value = 30241
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = ""

$ws.Range("H104").Value = 60000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 60000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 60000
$ws.Range("M104").Value = ""
$ws.Range("N104").Value = -65242

$ws.Range("H126").Value = 20001.2
$ws.Range("I126").Value = 110012
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 330036
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -327566
$ws.Range("N126").Value = -34940

$ws.Range("H132").Value = 63493724
$ws.Range("I132").Value = 63493724
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 190481172
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -190478642
$ws.Range("N132").Value = ""

$ws.Range("H141").Value = 128429.664
$ws.Range("J141").Value = 137784.44
$ws.Range("L141").Value = 137784.44
$ws.Range("N141").Value = -148144.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 48.875
$ws.Range("I2").Value = 46.333332
$ws.Range("J2").Value = 50.4
$ws.Range("K2").Value = 277.999992
$ws.Range("L2").Value = 302.4
$ws.Range("M2").Value = -164.999992
$ws.Range("N2").Value = -528.4

$ws.Range("H5").Value = 882.881
$ws.Range("I5").Value = 175.84616
$ws.Range("K5").Value = 527.5384799999999
$ws.Range("M5").Value = -415.5384799999999

$ws.Range("H23").Value = 510.7143
$ws.Range("J23").Value = 595.6
$ws.Range("L23").Value = 1786.8
$ws.Range("N23").Value = -2256.8

$ws.Range("H34").Value = 470446.9
$ws.Range("J34").Value = 2585.0667
$ws.Range("L34").Value = 7755.2001
$ws.Range("N34").Value = -7923.2001

$ws.Range("H75").Value = 1276.2727
$ws.Range("J75").Value = 1547.4286
$ws.Range("L75").Value = 4642.2858
$ws.Range("N75").Value = -6638.2858

$ws.Range("H78").Value = 1276.2727
$ws.Range("J78").Value = 1547.4286
$ws.Range("L78").Value = 13926.8574
$ws.Range("N78").Value = -23910.8574

$ws.Range("H121").Value = 755.1053000000001
$ws.Range("I121").Value = 257.33334
$ws.Range("J121").Value = 848.4375
$ws.Range("K121").Value = 772.0000200000001
$ws.Range("L121").Value = 2545.3125
$ws.Range("M121").Value = 537.9999799999999
$ws.Range("N121").Value = -5165.3125

$ws.Range("H135").Value = 882.881
$ws.Range("I135").Value = 175.84616
$ws.Range("K135").Value = 1582.61544
$ws.Range("M135").Value = 952.38456

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 28375
$ws.Range("J15").Value = 29500
$ws.Range("L15").Value = 29500
$ws.Range("N15").Value = -30076

$ws.Range("H81").Value = 28375
$ws.Range("J81").Value = 29500
$ws.Range("L81").Value = 29500
$ws.Range("N81").Value = -31496

$ws.Range("H84").Value = 28375
$ws.Range("J84").Value = 29500
$ws.Range("L84").Value = 88500
$ws.Range("N84").Value = -98484

$ws.Range("H97").Value = 1736
$ws.Range("I97").Value = 989.8
$ws.Range("J97").Value = 3601.5
$ws.Range("K97").Value = 989.8
$ws.Range("L97").Value = 3601.5
$ws.Range("M97").Value = -493.8
$ws.Range("N97").Value = -4593.5

$ws.Range("H122").Value = 361704.84
$ws.Range("I122").Value = 716287.7
$ws.Range("J122").Value = 7122
$ws.Range("K122").Value = 2148863.1
$ws.Range("L122").Value = 21366
$ws.Range("M122").Value = -2146413.1
$ws.Range("N122").Value = -26266

$ws.Range("H123").Value = 45289.5
$ws.Range("J123").Value = 45289.5
$ws.Range("L123").Value = 45289.5
$ws.Range("N123").Value = -50189.5

$ws.Range("H126").Value = 3327.524
$ws.Range("I126").Value = 2400
$ws.Range("K126").Value = 7200
$ws.Range("M126").Value = -4730

$ws.Range("H133").Value = 125000
$ws.Range("J133").Value = 125000
$ws.Range("L133").Value = 125000
$ws.Range("N133").Value = -135120

$ws.Range("H141").Value = 31843
$ws.Range("J141").Value = 31843
$ws.Range("L141").Value = 31843
$ws.Range("N141").Value = -42203

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 28649666
$ws.Range("I40").Value = 13891412
$ws.Range("K40").Value = 13891412
$ws.Range("M40").Value = -13891276

$ws.Range("H42").Value = 31666.666
$ws.Range("I42").Value = 30000
$ws.Range("K42").Value = 30000
$ws.Range("M42").Value = -29437

$ws.Range("H49").Value = 31666.666
$ws.Range("I49").Value = 30000
$ws.Range("K49").Value = 30000
$ws.Range("M49").Value = -29853

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""

$ws.Range("H136").Value = 3952.4807
$ws.Range("I136").Value = 2551.647
$ws.Range("K136").Value = 7654.941
$ws.Range("M136").Value = -5104.941

$ws.Range("H138").Value = 96098
$ws.Range("I138").Value = 94900
$ws.Range("J138").Value = 96697
$ws.Range("K138").Value = 94900
$ws.Range("L138").Value = 96697
$ws.Range("M138").Value = -89760
$ws.Range("N138").Value = -106977

$ws.Range("H140").Value = 96614
$ws.Range("J140").Value = 96614
$ws.Range("L140").Value = 96614
